$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "Problems" column (E) for chapter 10 and chapter 11 rows.
$ws.Range("E11").Value = 12
$ws.Range("E12").Value = 5

# Move the active selection from F13 to E13 (G3 formula recalculates automatically).
$ws.Range("E13").Select()
